$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '306.01'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.91%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '39.23'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '7.84%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.116'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.15%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08070'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.938'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-2.21%'

$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.199'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.15%'

$ws.Range("B8").Value = 'KuCoinToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '8.060'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2.51%'

$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9255'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.07%'

$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1382'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-6.74%'

$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1918'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.72%'

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08993'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.87%'

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03520'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.12%'

$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09792'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.68%'

$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001398'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.19%'

$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005932'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-6.02%'

$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.770'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-2.00%'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.25%'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3462'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1294'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.44%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.673'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.81%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2417'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '3.14%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04367'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.15%'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001207'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-2.22%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004281'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '2.85%'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.05%'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02037'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-2.48%'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.78%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007520'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.56%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009782'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-3.09%'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.73%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002093'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-1.37%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009801'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.07%'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006209'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-1.15%'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.10%'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002779'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '12.58%'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.10%'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.10%'
